# Atualização automática de PANAMBI.xlsx
#
# - Renomeia "Paineis DARQ"            -> "PAINEIS DARQ"
# - Renomeia "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# - Remove a planilha "Desarquivamentos Pendentes"

$wb = $excel.ActiveWorkbook

# Avoid the "delete sheet" confirmation dialog.
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

[void]$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

$excel.DisplayAlerts = $true
